$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.80779999999999
$ws.Range("E4").Value = 12.74200000000001
$ws.Range("B9").Value = 8.833200000000007
$ws.Range("E10").Value = 12.13349999999999
$ws.Range("B18").Value = 4.590400000000005
$ws.Range("B20").Value = 5.800800000000001
$ws.Range("D21").Value = -7.325000000000005
